$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.906.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.814.84'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3661'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07354'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8685'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.884.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.385'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07092'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.515'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008702'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.931.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.060.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.894'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.150'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.265'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08897'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7542'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.157'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.490'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.912'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05274'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01948'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.977'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.244'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5305'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.291'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.90%  '
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.430'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("E46").Value = '  -2.69%  '
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.21'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.660'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06292'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.13%  '
